# Auto-generated Excel COM-interop edit script
# Updates leve-profit market-price columns (H-N) across all 8 job sheets
# per the scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 8
$ws.Range("H8").Value = 339.91666
$ws.Range("I8").Value = 339.91666
$ws.Range("K8").Value = 1019.74998
$ws.Range("M8").Value = -880.7499799999999

# Row 113
$ws.Range("H113").Value = 4388.737
$ws.Range("I113").Value = 5600
$ws.Range("J113").Value = 3682.1667
$ws.Range("K113").Value = 5600
$ws.Range("L113").Value = 3682.1667
$ws.Range("M113").Value = -2346
$ws.Range("N113").Value = -10190.1667

# Row 116
$ws.Range("H116").Value = 5668712.5
$ws.Range("J116").Value = 2285.5557
$ws.Range("L116").Value = 2285.5557
$ws.Range("N116").Value = -9169.555700000001

# Row 132
$ws.Range("H132").Value = 5490.1763
$ws.Range("I132").Value = 7047.722
$ws.Range("J132").Value = 3737.9375
$ws.Range("K132").Value = 21143.166
$ws.Range("L132").Value = 11213.8125
$ws.Range("M132").Value = -18613.166
$ws.Range("N132").Value = -16273.8125

# Row 138
$ws.Range("H138").Value = 2114.0232
$ws.Range("I138").Value = 702.61536
$ws.Range("J138").Value = 4272.647
$ws.Range("K138").Value = 2107.84608
$ws.Range("L138").Value = 12817.941
$ws.Range("M138").Value = 3032.15392
$ws.Range("N138").Value = -23097.941

$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 9262711
$ws.Range("I74").Value = 16670380
$ws.Range("J74").Value = 3125.1667
$ws.Range("K74").Value = 16670380
$ws.Range("L74").Value = 3125.1667
$ws.Range("M74").Value = -16669506
$ws.Range("N74").Value = -4873.1667

# Row 77
$ws.Range("H77").Value = 9262711
$ws.Range("I77").Value = 16670380
$ws.Range("J77").Value = 3125.1667
$ws.Range("K77").Value = 83351900
$ws.Range("L77").Value = 15625.8335
$ws.Range("M77").Value = -83347532
$ws.Range("N77").Value = -24361.8335

# Row 88
$ws.Range("H88").Value = 90911620
$ws.Range("I88").Value = 2078
$ws.Range("J88").Value = 166669570
$ws.Range("K88").Value = 2078
$ws.Range("L88").Value = 166669570
$ws.Range("M88").Value = -1672
$ws.Range("N88").Value = -166670382

# Row 91
$ws.Range("H91").Value = 90911620
$ws.Range("I91").Value = 2078
$ws.Range("J91").Value = 166669570
$ws.Range("K91").Value = 2078
$ws.Range("L91").Value = 166669570
$ws.Range("M91").Value = -674
$ws.Range("N91").Value = -166672378

# Row 97
$ws.Range("H97").Value = 1796.8695
$ws.Range("I97").Value = 2253.353
$ws.Range("J97").Value = 503.5
$ws.Range("K97").Value = 2253.353
$ws.Range("L97").Value = 503.5
$ws.Range("M97").Value = -1757.353
$ws.Range("N97").Value = -1495.5

# Row 132
$ws.Range("H132").Value = 2044.8077
$ws.Range("I132").Value = 1795.9556
$ws.Range("J132").Value = 3644.5715
$ws.Range("K132").Value = 5387.8668
$ws.Range("L132").Value = 10933.7145
$ws.Range("M132").Value = -2857.8668
$ws.Range("N132").Value = -15993.7145

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 6458.8335
$ws.Range("I20").Value = 3598.75
$ws.Range("J20").Value = 7888.875
$ws.Range("K20").Value = 3598.75
$ws.Range("L20").Value = 7888.875
$ws.Range("M20").Value = -3351.75
$ws.Range("N20").Value = -8382.875

# Row 86
$ws.Range("H86").Value = 10001671
$ws.Range("I86").Value = 12501589
$ws.Range("K86").Value = 12501589
$ws.Range("M86").Value = -12500466

# Row 89
$ws.Range("H89").Value = 10001671
$ws.Range("I89").Value = 12501589
$ws.Range("K89").Value = 62507945
$ws.Range("M89").Value = -62502329

# Row 94
$ws.Range("H94").Value = 13896.6
$ws.Range("I94").Value = 492.07693
$ws.Range("J94").Value = 101026
$ws.Range("K94").Value = 492.07693
$ws.Range("L94").Value = 101026
$ws.Range("M94").Value = -41.07693
$ws.Range("N94").Value = -101928

# Row 134
$ws.Range("H134").Value = 2090.0435
$ws.Range("I134").Value = 1653
$ws.Range("J134").Value = 3328.3333
$ws.Range("K134").Value = 4959
$ws.Range("L134").Value = 9984.999899999999
$ws.Range("M134").Value = -2424
$ws.Range("N134").Value = -15054.9999

$ws = $wb.Worksheets.Item("CRP")
# Row 9
$ws.Range("H9").Value = 97500
$ws.Range("J9").Value = 97500
$ws.Range("L9").Value = 97500
$ws.Range("N9").Value = -97836

# Row 105
$ws.Range("H105").Value = 829.6667
$ws.Range("I105").Value = 829.6667
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 829.6667
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 917.3333
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 1000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 3000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3338

# Row 34
$ws.Range("H34").Value = 472.27585
$ws.Range("I34").Value = 59.6
$ws.Range("J34").Value = 689.4737
$ws.Range("K34").Value = 178.8
$ws.Range("L34").Value = 2068.4211
$ws.Range("M34").Value = -94.80000000000001
$ws.Range("N34").Value = -2236.4211

# Row 39
$ws.Range("H39").Value = 1760
$ws.Range("I39").Value = 1000
$ws.Range("J39").Value = 1950
$ws.Range("K39").Value = 3000
$ws.Range("L39").Value = 5850
$ws.Range("M39").Value = -2706
$ws.Range("N39").Value = -6438

# Row 55
$ws.Range("H55").Value = 1564.3334
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1564.3334
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 4693.0002
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -5047.0002

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 33338598
$ws.Range("I70").Value = 80004780
$ws.Range("K70").Value = 80004780
$ws.Range("M70").Value = -80004510

# Row 73
$ws.Range("H73").Value = 33338598
$ws.Range("I73").Value = 80004780
$ws.Range("K73").Value = 80004780
$ws.Range("M73").Value = -80003844

# Row 97
$ws.Range("H97").Value = 2503.3333
$ws.Range("I97").Value = 2503.3333
$ws.Range("K97").Value = 2503.3333
$ws.Range("M97").Value = -2007.3333

# Row 113
$ws.Range("H113").Value = 1672.1666
$ws.Range("I113").Value = 1506.6
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 1506.6
$ws.Range("L113").Value = 2500
$ws.Range("M113").Value = 663.4000000000001
$ws.Range("N113").Value = -6840

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2223308
$ws.Range("I22").Value = 3334143.2
$ws.Range("J22").Value = 1637.6
$ws.Range("K22").Value = 3334143.2
$ws.Range("L22").Value = 1637.6
$ws.Range("M22").Value = -3333848.2
$ws.Range("N22").Value = -2227.6

# Row 27
$ws.Range("H27").Value = 2223308
$ws.Range("I27").Value = 3334143.2
$ws.Range("J27").Value = 1637.6
$ws.Range("K27").Value = 3334143.2
$ws.Range("L27").Value = 1637.6
$ws.Range("M27").Value = -3334036.2
$ws.Range("N27").Value = -1851.6

# Row 61
$ws.Range("H61").Value = 1170.4667
$ws.Range("I61").Value = 1242.6666
$ws.Range("J61").Value = 881.6667
$ws.Range("K61").Value = 1242.6666
$ws.Range("L61").Value = 881.6667
$ws.Range("M61").Value = -1040.6666
$ws.Range("N61").Value = -1285.6667

# Row 113
$ws.Range("H113").Value = 1170.4667
$ws.Range("I113").Value = 1242.6666
$ws.Range("J113").Value = 881.6667
$ws.Range("K113").Value = 1242.6666
$ws.Range("L113").Value = 881.6667
$ws.Range("M113").Value = 927.3334
$ws.Range("N113").Value = -5221.6667

# Row 122
$ws.Range("H122").Value = 6251
$ws.Range("I122").Value = 5752
$ws.Range("K122").Value = 17256
$ws.Range("M122").Value = -14806

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 1266.6666
$ws.Range("I96").Value = 1266.6666
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = 106.3334
$ws.Range("N96").ClearContents()
